$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

$ws.Range("M1").Value = "portraitAddress|String"
$ws.Range("M2").Value = "Portrati_Ganfaul"
$ws.Range("M3").Value = "Portrati_KeepSeries"
$ws.Range("M4").Value = "Portrati_BigBatSuccubus"
$ws.Range("M5").Value = "Portrati_Bei"

$ws.Columns.Item(13).ColumnWidth = 16.85
